# Update image link paths (column B / Bild_Link) for rows 4-11 on the "Laptops" sheet
# so that each laptop model links to its own matching image instead of the
# previously mismatched "Dell Precision 5550 ..." / "Dell Latitude 3550 ..." images.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Laptops")

$ws.Range("B4").Value = "images/Dell Precision 3561 i9 11950H 15 Zoll Notebook Workstation in Brandenburg - Frankfurt (Oder) _ kleinanzeigen.de/Dell Precision 3561 i9 11950H 15 Zoll Notebook Workstation in Brandenburg - Frankfurt (Oder) _ kleinanzeigen.de_page2_img1.png"
$ws.Range("B5").Value = "images/Dell Precision 7560 _ i7-11850H _ 15.6_ Refurbished/Dell Precision 7560 _ i7-11850H _ 15.6_ Refurbished_page1_img1.png"
$ws.Range("B6").Value = "images/Dell Precision 7550 - i7-10750H - 15.6_ Laptop2/Dell Precision 7550 - i7-10750H - 15.6_ Laptop2_page1_img1.png"
$ws.Range("B7").Value = "images/Dell Latitude 5501 _ i7-9850H _ 15.6_ - Refurbished/Dell Latitude 5501 _ i7-9850H _ 15.6_ - Refurbished_page1_img1.png"
$ws.Range("B8").Value = "images/Dell Precision 7540 _ i9-9880H _ 15.6_ - Refurbished/Dell Precision 7540 _ i9-9880H _ 15.6_ - Refurbished_page1_img1.png"
$ws.Range("B9").Value = "images/Dell Precision 5560 _ i5-11500H _ 15.6'' - refurbished/Dell Precision 5560 _ i5-11500H _ 15.6'' - refurbished_page1_img1.png"
$ws.Range("B10").Value = "images/Dell Precision 5560 _ i5-11500H _ 15.6'' - refurbished/Dell Precision 5560 _ i5-11500H _ 15.6'' - refurbished_page1_img1.png"
$ws.Range("B11").Value = "images/Dell Precision 5560 _ i5-11500H _ 15.6'' - refurbished/Dell Precision 5560 _ i5-11500H _ 15.6'' - refurbished_page1_img1.png"
